$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$naText = "N/A - Stat tracked as of the 1973-74 ABA Season"

# stl_per_game (AQ) and blk_per_game (AR) were not tracked before the
# 1973-74 ABA season. Replace the stray/missing values in those rows with
# an explanatory N/A note.
$ws.Range("AQ6").Value = $naText
$ws.Range("AR6").Value = $naText

$ws.Range("AQ7").Value = $naText
$ws.Range("AR7").Value = $naText

$ws.Range("AQ8").Value = $naText
$ws.Range("AR8").Value = $naText

$ws.Range("AQ9").Value = $naText
$ws.Range("AR9").Value = $naText

$ws.Range("AQ10").Value = $naText
$ws.Range("AR10").Value = $naText

$ws.Range("AQ11").Value = $naText
$ws.Range("AR11").Value = $naText
